$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.381.36'
$ws.Range("E2").Value = '  +1.82%  '
$ws.Range("D3").Value = '1.883.14'
$ws.Range("E3").Value = '  +0.43%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '0.700'
$ws.Range("E5").Value = '  +1.69%  '
$ws.Range("D6").Value = '246.87'
$ws.Range("E6").Value = '  -0.49%  '
$ws.Range("D8").Value = '43.35'
$ws.Range("E8").Value = '  +4.79%  '
$ws.Range("D9").Value = '0.357'
$ws.Range("E9").Value = '  +2.71%  '
$ws.Range("D10").Value = '0.0747'
$ws.Range("E10").Value = '  +1.03%  '
$ws.Range("E11").Value = '  +1.18%  '
$ws.Range("D12").Value = '13.51'
$ws.Range("E12").Value = '  +5.04%  '
$ws.Range("D13").Value = '2.156.85'
$ws.Range("E13").Value = '  +0.43%  '
$ws.Range("D14").Value = '0.773'
$ws.Range("E14").Value = '  +8.25%  '
$ws.Range("E15").Value = '  +1.11%  '
$ws.Range("D16").Value = '1.857.14'
$ws.Range("E16").Value = '  -0.75%  '
$ws.Range("D17").Value = '35.357.89'
$ws.Range("E17").Value = '  +1.79%  '
$ws.Range("D18").Value = '73.45'
$ws.Range("E18").Value = '  +0.81%  '
$ws.Range("D19").Value = '0.0₃0828'
$ws.Range("E19").Value = '  +0.55%  '
$ws.Range("D20").Value = '244.69'
$ws.Range("E20").Value = '  -1.17%  '
$ws.Range("D21").Value = '12.84'
$ws.Range("E21").Value = '  +0.61%  '
$ws.Range("D22").Value = '5.18'
$ws.Range("E22").Value = '  +5.21%  '
$ws.Range("E23").Value = '  +7.82%  '
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("D25").Value = '2.16'
$ws.Range("E25").Value = '  -3.93%  '
$ws.Range("D26").Value = '164.65'
$ws.Range("E26").Value = '  -0.40%  '
$ws.Range("E27").Value = '  +3.02%  '
$ws.Range("D28").Value = '18.32'
$ws.Range("E28").Value = '  +0.38%  '
$ws.Range("D29").Value = '0.129'
$ws.Range("E29").Value = '  +0.79%  '
$ws.Range("D30").Value = '0.0598'
$ws.Range("E30").Value = '  +3.14%  '
$ws.Range("D31").Value = '4.29'
$ws.Range("E31").Value = '  +0.95%  '
$ws.Range("E32").Value = '  +3.51%  '
$ws.Range("D33").Value = '4.19'
$ws.Range("E33").Value = '  +0.47%  '
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("E35").Value = '  -12.51%  '
$ws.Range("D36").Value = '0.855'
$ws.Range("E36").Value = '  +2.26%  '
$ws.Range("E37").Value = '  -1.17%  '
$ws.Range("E38").Value = '  +11.67%  '
$ws.Range("D39").Value = '17.31'
$ws.Range("E39").Value = '  +0.28%  '
$ws.Range("E40").Value = '  +3.69%  '
$ws.Range("D41").Value = '97.30'
$ws.Range("E41").Value = '  -1.03%  '
$ws.Range("E42").Value = '  -0.69%  '
$ws.Range("E43").Value = '  +1.94%  '
$ws.Range("D44").Value = '1.308.23'
$ws.Range("E44").Value = '  +1.13%  '
$ws.Range("E45").Value = '  +6.08%  '
$ws.Range("E46").Value = '  -0.54%  '
$ws.Range("D47").Value = '2.74'
$ws.Range("E47").Value = '  +0.52%  '
$ws.Range("E48").Value = '  +0.33%  '
$ws.Range("E49").Value = '  -2.53%  '
$ws.Range("D50").Value = '42.35'
$ws.Range("E50").Value = '  +0.61%  '
$ws.Range("D51").Value = '2.060.50'
$ws.Range("E51").Value = '  +0.16%  '
